$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The figures in this table are stored as text (not numbers), so we
# temporarily force a text NumberFormat while writing the new values to
# keep Excel from auto-converting them to numeric cells, then restore the
# original "Normal" cell style so no stray formatting is left behind.

# Enterprises density (per 1000 people) - row 13
$r1 = $ws.Range("B13:D13")
$r1.NumberFormat = "@"
$ws.Range("B13").Value = "12.86"
$ws.Range("C13").Value = "0.65"
$ws.Range("D13").Value = "13.52"
$r1.Style = "Normal"

# Employment (% of total) - row 14
$r2 = $ws.Range("B14:D14")
$r2.NumberFormat = "@"
$ws.Range("B14").Value = "67.96"
$ws.Range("C14").Value = "18.66"
$ws.Range("D14").Value = "86.61"
$r2.Style = "Normal"

# Enterprises (% of total) - row 16
$r3 = $ws.Range("B16:D16")
$r3.NumberFormat = "@"
$ws.Range("B16").Value = "94.45"
$ws.Range("C16").Value = "4.79"
$ws.Range("D16").Value = "99.24"
$r3.Style = "Normal"
